$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Date and Time'
$ws.Range("B1").Style = "Normal"
$ws.Range("B1").Value = "2024-03-11 11:56:50.902000 to 2024-03-11 12:41:15.782000"

$ws.Range("A2").Value = 'Total time taken for the ride'
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"
$ws.Range("B2").Value = 0.03086597222222223

$ws.Range("A3").Value = 'Actual Ampere-hours (Ah)'
$ws.Range("B3").Value = 27.4824425

$ws.Range("A4").Value = 'Actual Watt-hours (Wh)'
$ws.Range("B4").Value = 1392.512839500555

$ws.Range("A5").Value = 'Starting SoC (Ah)'
$ws.Range("B5").Value = 38.734

$ws.Range("A6").Value = 'Ending SoC (Ah)'
$ws.Range("B6").Value = 11.11

$ws.Range("A7").Value = 'Starting SoC (%)'
$ws.Range("B7").Value = 28

$ws.Range("A8").Value = 'Ending SoC (%)'
$ws.Range("B8").Value = 97

$ws.Range("A9").Value = 'Total distance covered (km)'
$ws.Range("B9").Value = 28.99043299372953

$ws.Range("A10").Value = 'Total energy consumption(WH/KM)'
$ws.Range("B10").Value = 48.03353022708382

$ws.Range("A11").Value = 'Total SOC consumed(%)'
$ws.Range("B11").Value = 69

$ws.Range("A12").Value = 'Mode'
$ws.Range("B12").Value = "Custom mode`n94.60%`nEco mode`n3.25%`nSports mode`n0.11%"

$ws.Range("A13").Value = 'Peak Power(kW)'
$ws.Range("B13").Value = 5527.275720000001

$ws.Range("A14").Value = 'Average Power(kW)'
$ws.Range("B14").Value = -1882.48074434923

$ws.Range("A15").Value = 'Total Energy Regenerated(kWh)'
$ws.Range("B15").Value = 0.1260407530555556

$ws.Range("A16").Value = 'Regenerative Effectiveness(%)'
$ws.Range("B16").Value = 0.009050497931854558

$ws.Range("A17").Value = 'Highest Cell Voltage(V)'
$ws.Range("B17").Value = 3.328

$ws.Range("A18").Value = 'Lowest Cell Voltage(V)'
$ws.Range("B18").Value = 3.06

$ws.Range("A19").Value = 'Difference in Cell Voltage(V)'
$ws.Range("B19").Value = 0.2679999999999998

$ws.Range("A20").Value = 'Minimum Temperature(C)'
$ws.Range("B20").Value = 28

$ws.Range("A21").Value = 'Maximum Temperature(C)'
$ws.Range("B21").Value = 45

$ws.Range("A22").Value = 'Difference in Temperature(C)'
$ws.Range("B22").Value = 17

$ws.Range("A23").Value = 'Maximum Fet Temperature-BMS(C)'
$ws.Range("B23").Value = 70

$ws.Range("A24").Value = 'Maximum Afe Temperature-BMS(C)'
$ws.Range("B24").Value = 67

$ws.Range("A25").Value = 'Maximum PCB Temperature-BMS(C)'
$ws.Range("B25").Value = 66

$ws.Range("A26").Value = 'Maximum MCU Temperature(C)'
$ws.Range("B26").Value = 64

$ws.Range("A27").Value = 'Maximum Motor Temperature(C)'
$ws.Range("B27").Value = 102

$ws.Range("A28").Value = 'Abnormal Motor Temperature Detected(C)'
$ws.Range("B28").Value = 0

$ws.Range("A29").Value = 'highest cell temp(C)'
$ws.Range("B29").Value = 45

$ws.Range("A30").Value = 'lowest cell temp(C)'
$ws.Range("B30").Value = 28

$ws.Range("A31").Value = 'Difference between Highest and Lowest Cell Temperature at 100% SOC(C)'
$ws.Range("B31").Value = 17

$ws.Range("A32").Value = 'Battery Voltage(V)'
$ws.Range("B32").Value = 53

$ws.Range("A33").Value = 'Total energy charged(kWh)'
$ws.Range("B33").Value = 1.4565694525

$ws.Range("A34").Value = 'Electricity consumption units(kW)'
$ws.Range("B34").Value = 0.0000001517639256168209

$ws.Range("A35").Value = 'Cycle Count of battery'
$ws.Range("B35").Value = 39

$ws.Range("A36").Value = 'Idling time percentage'
$ws.Range("B36").Value = 14.16272583969928

$ws.Range("A37").Value = 'Time spent in 0-10 km/h'
$ws.Range("B37").Value = 5.363566549452326

$ws.Range("A38").Value = 'Time spent in 10-20 km/h'
$ws.Range("B38").Value = 6.798431752960672

$ws.Range("A39").Value = 'Time spent in 20-30 km/h'
$ws.Range("B39").Value = 10.02384705549493

$ws.Range("A40").Value = 'Time spent in 30-40 km/h'
$ws.Range("B40").Value = 10.97772927529203

$ws.Range("A41").Value = 'Time spent in 40-50 km/h'
$ws.Range("B41").Value = 9.166969807202619

$ws.Range("A42").Value = 'Time spent in 50-60 km/h'
$ws.Range("B42").Value = 15.56929792651874

$ws.Range("A43").Value = 'Time spent in 60-70 km/h'
$ws.Range("B43").Value = 15.03981245705509

$ws.Range("A44").Value = 'Time spent in 70-80 km/h'
$ws.Range("B44").Value = 12.76423750050524

$ws.Range("A45").Value = 'Time spent in 80-90 km/h'
$ws.Range("B45").Value = 0
